$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-4 ---
$ws.Range("F2").Value = "2G_BVI010M_HNI"
$ws.Range("G2").Value = "Thon-Lien-Tong-BVI_HNI"
$ws.Range("J2").Value = "POWER_AC_EAS"
$ws.Range("L2").Value = "06/05/2025 15:42:22"
$ws.Range("T2").Value = "Ba Vì"
$ws.Range("AA2").Value = "Trạm viễn thông loại 1"

$ws.Range("F3").Value = "3G_BVI010M_HNI"
$ws.Range("G3").Value = "TONG-BAT-THON-TONG-LENH-BVI_HNI"
$ws.Range("J3").Value = "POWER_AC_EAS"
$ws.Range("L3").Value = "06/05/2025 15:37:25"
$ws.Range("T3").Value = "Ba Vì"
$ws.Range("V3").Value = ""
$ws.Range("AA3").Value = "Trạm viễn thông loại 3"

$ws.Range("F4").Value = "4G-BVI010M-HNI"
$ws.Range("G4").Value = "KCN-CAU-GAO-DPG_HNI"
$ws.Range("J4").Value = "POWER_AC_EAS"
$ws.Range("L4").Value = "06/05/2025 14:30:14"
$ws.Range("T4").Value = "Đan Phượng"
$ws.Range("AA4").Value = "Trạm viễn thông loại 3"

# --- Add new rows 5-8, copying style from row 4 ---
$ws.Range("A4:AC4").Copy()
$ws.Range("A5:AC8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Clear previous content/formula, keep just the pasted style, then set values
$ws.Range("A5:AC8").ClearContents()

# Row 5
$ws.Range("F5").Value = "2G_BVI010M_HNI"
$ws.Range("G5").Value = "Cam-Thuong-Thon-Van-Minh-BVI_HNI"
$ws.Range("J5").Value = "SITE_OOS_BY_POWER"
$ws.Range("L5").Value = "06/05/2025 13:42:12"
$ws.Range("T5").Value = "Ba Vì"
$ws.Range("V5").Value = "184205 - VTHN TĐML - HNI dựng lại côt thay cáp - 4 - thainh1 - 06/05/2025 13:56:39"
$ws.Range("AA5").Value = "Trạm viễn thông loại 2"

# Row 6
$ws.Range("F6").Value = "3G_BVI010M_HNI"
$ws.Range("G6").Value = "Cam-Thuong-Thon-Van-Minh-BVI_HNI"
$ws.Range("J6").Value = "SITE_OOS_BY_POWER"
$ws.Range("L6").Value = "06/05/2025 13:41:56"
$ws.Range("T6").Value = "Ba Vì"
$ws.Range("V6").Value = "184205 - VTHN TĐML - HNI dựng lại côt thay cáp - 4 - thainh1 - 06/05/2025 13:56:40"
$ws.Range("AA6").Value = "Trạm viễn thông loại 2"

# Row 7
$ws.Range("F7").Value = "4G-BVI010M-HNI"
$ws.Range("G7").Value = "Cam-Thuong-Thon-Van-Minh-BVI_HNI"
$ws.Range("J7").Value = "SITE_OOS_BY_POWER"
$ws.Range("L7").Value = "06/05/2025 13:41:51"
$ws.Range("T7").Value = "Ba Vì"
$ws.Range("V7").Value = "184205 - VTHN TĐML - HNI dựng lại côt thay cáp - 4 - thainh1 - 06/05/2025 13:56:40"
$ws.Range("AA7").Value = "Trạm viễn thông loại 2"

# Row 8
$ws.Range("F8").Value = "SR_BVI010M_HNI"
$ws.Range("G8").Value = "Cam-Thuong-Thon-Van-Minh-BVI_HNI"
$ws.Range("J8").Value = "POWER_AC_EAS"
$ws.Range("L8").Value = "06/05/2025 06:18:28"
$ws.Range("T8").Value = "Ba Vì"
$ws.Range("V8").Value = "Mất nguồn AC - 1 - huongvl1 - 06/05/2025 10:23:32"
$ws.Range("AA8").Value = "Trạm viễn thông loại 2"

# --- Column width adjustments ---
# (Target raw widths are 34.7109375 / 19.7109375 / 84.7109375 "character" units;
#  the ColumnWidth setter here snaps to a 1/6 grid offset by 5/6, so we pick the
#  nearest representable input to land as close as possible to the target.)
$ws.Columns.Item(7).ColumnWidth = 33.833333333333336
$ws.Columns.Item(10).ColumnWidth = 18.833333333333332
$ws.Columns.Item(22).ColumnWidth = 83.83333333333333
